$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 31 : vial 30 dies (censor=1) on day 48, date 45632 -> style = red
#          (same "dead" style already used elsewhere, e.g. B2:C2)
# ---------------------------------------------------------------------
$ws.Range("B2:C2").Copy()
$ws.Range("B31:C31").PasteSpecial(-4122)
$ws.Cells.Item(31,3).Value = 48

$ws.Cells.Item(2,4).Copy()
$ws.Cells.Item(31,4).PasteSpecial(-4122)
$ws.Cells.Item(31,4).Value = 45632

$ws.Cells.Item(31,5).Value = 1

# ---------------------------------------------------------------------
# Row 33 : vial 32 dies (censor=1) on day 50, date 45634 -> new style
#          (red fill, same as the "dead" style but also flags the font)
# ---------------------------------------------------------------------
$ws.Cells.Item(33,2).Interior.Color = 255
$ws.Cells.Item(33,2).Font.ThemeColor = 1
$ws.Cells.Item(33,2).Font.TintAndShade = 0

$ws.Cells.Item(33,2).Copy()
$ws.Cells.Item(33,3).PasteSpecial(-4122)
$ws.Cells.Item(33,3).Value = 50

$ws.Cells.Item(2,4).Copy()
$ws.Cells.Item(33,4).PasteSpecial(-4122)
$ws.Cells.Item(33,4).Value = 45634

$ws.Cells.Item(33,5).Value = 1

# ---------------------------------------------------------------------
# Row 36 : vial 35 still alive (censor=0) at day 48, date 45632
#          -> blue "still alive" style already used elsewhere (B21:C21)
# ---------------------------------------------------------------------
$ws.Range("B21:C21").Copy()
$ws.Range("B36:C36").PasteSpecial(-4122)
$ws.Cells.Item(36,3).Value = 48

$ws.Cells.Item(2,4).Copy()
$ws.Cells.Item(36,4).PasteSpecial(-4122)
$ws.Cells.Item(36,4).Value = 45632

$ws.Cells.Item(36,5).Value = 0

# ---------------------------------------------------------------------
# Row 78 : vial 77 dies (censor=1) on day 50, date 45634 -> red style
# ---------------------------------------------------------------------
$ws.Range("B2:C2").Copy()
$ws.Range("B78:C78").PasteSpecial(-4122)
$ws.Cells.Item(78,3).Value = 50

$ws.Cells.Item(2,4).Copy()
$ws.Cells.Item(78,4).PasteSpecial(-4122)
$ws.Cells.Item(78,4).Value = 45634

$ws.Cells.Item(78,5).Value = 1

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Update the view: scroll so row 27 is at the top, and select D33
# ---------------------------------------------------------------------
$ws.Application.Goto($ws.Range("A27"))
$ws.Range("D33").Select()
